$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "47.732.39"
$ws.Range("E2").Value = "  +5.55%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.506.69"
$ws.Range("E3").Value = "  +3.14%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "323.96"
$ws.Range("E5").Value = "  +2.19%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "106.03"
$ws.Range("E6").Value = "  +3.22%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.525"
$ws.Range("E7").Value = "  +1.63%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("E9").Value = "  +2.93%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.07"
$ws.Range("E10").Value = "  +7.13%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0818"
$ws.Range("E11").Value = "  +1.80%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.123"
$ws.Range("E12").Value = "  +0.93%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.43"
$ws.Range("E13").Value = "  +1.65%  "
$ws.Range("E14").Value = "  +2.22%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.901.89"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.507.82"
$ws.Range("E16").Value = "  +3.30%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.848"
$ws.Range("E17").Value = "  +0.85%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "47.624.68"
$ws.Range("E18").Value = "  +5.59%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.81"
$ws.Range("E19").Value = "  +4.40%  "
$ws.Range("E20").Value = "  +3.47%  "
$ws.Range("E21").Value = "  +1.97%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "70.82"
$ws.Range("E22").Value = "  +2.86%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "251.45"
$ws.Range("E23").Value = "  +3.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.39"
$ws.Range("E24").Value = "  +6.16%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.58"
$ws.Range("E25").Value = "  +3.21%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.32"
$ws.Range("E26").Value = "  +2.80%  "
$ws.Range("E27").Value = "  -0.10%  "
$ws.Range("B28").Value = "Cosmos"
$ws.Range("C28").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.04"
$ws.Range("E28").Value = "  +4.86%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.21"
$ws.Range("E29").Value = "  +6.70%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.19"
$ws.Range("E30").Value = "  +6.52%  "
$ws.Range("E31").Value = "  +8.91%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.47"
$ws.Range("E32").Value = "  +0.60%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.99"
$ws.Range("E33").Value = "  -1.72%  "
$ws.Range("E34").Value = "  +3.15%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0786"
$ws.Range("E35").Value = "  +2.85%  "
$ws.Range("E36").Value = "  +0.13%  "
$ws.Range("E37").Value = "  +3.99%  "
$ws.Range("E38").Value = "  +4.50%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.00"
$ws.Range("E39").Value = "  +5.05%  "
$ws.Range("B40").Value = "Stellar"
$ws.Range("C40").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.112"
$ws.Range("E40").Value = "  +2.17%  "
$ws.Range("B41").Value = "WEMIXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.26"
$ws.Range("E41").Value = "  +2.19%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "121.65"
$ws.Range("E42").Value = "  -1.67%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.35"
$ws.Range("E43").Value = "  +2.96%  "
$ws.Range("E44").Value = "  +3.56%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.971.31"
$ws.Range("E45").Value = "  +1.81%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.02"
$ws.Range("E46").Value = "  +3.07%  "
$ws.Range("E47").Value = "  -0.38%  "
$ws.Range("E48").Value = "  +0.44%  "
$ws.Range("E49").Value = "  +0.72%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.36"
$ws.Range("E50").Value = "  +13.51%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "79.29"
$ws.Range("E51").Value = "  +3.51%  "
